$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture references to the rows that will shift down when we insert,
# so we can restore/copy their original formatting after the structural change.
$origFirstBlankRow = $ws.Range("A34:G34")   # currently styled as the "top of blank block" row
$origDataRowFormat = $ws.Range("A33:G33")   # last existing diary entry row (donor format for new entries)

# Insert 5 new blank rows at 34:38 to make room for the new diary entries;
# this pushes the existing blank rows (34-39) down to (39-44).
$ws.Rows.Item("34:38").Insert()

# Re-apply the original "top of blank block" formatting to its new location (row 39).
$origFirstBlankRow.Copy()
$ws.Range("A39:G39").PasteSpecial(-4122)

# Copy the diary-entry formatting onto the 5 newly inserted rows.
$origDataRowFormat.Copy()
$ws.Range("A34:G38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 34 (17:00 - 19:50) ---
$ws.Rows.Item(34).RowHeight = 376.7
$ws.Range("A34").Value = 43881
$ws.Range("B34").Value = "17:00 - 19:50"
$ws.Range("C34").Value = "None"
$ws.Range("D34").Value = "Looking forward to discuss the assignment as we had hard time settling on the essential functional and non-functional features, want to hear other groups’ opinions. Let’s see if we get our midterms today, kind of nervous. Also, looking forward to the speaker!"
$ws.Range("E34").Value = "We discussed on the features we found interesting, it was nice to hear the different security concerns or portability features from others’ systems. We studied few more key concepts. Followed by crash coursing Software Architecture. We motivated the architecture behind jPacman 3 with partners. We also learnt about the social context and the different standards that one should adhere to while making a pull request. Finally, had an interesting conversation with Sara and Omar!"
$ws.Range("F34").Value = "When questioned about the suppliers, nobody that was questioned was able to answer the details which was sad since we spent the previous week working for the assignment. While looking for the architectural pattern behind jPacman, I understood the various ways of thinking, bottom up or higher level abstractions.  Even though, we all knew roughly about the general principles, doing the exercise made as not so confident about our understanding of MVC and the arrows that connect these sectors. It was easier to understand once we grouped different sections in UML. Finally, it was interesting to listen to Sara and Omar. Sara seemed nervous in the newer settings which again reassures that even if we are experts in our field, it is normal to feel human, to be tensed when addressing a crowd, etc. Personally, my favorite till date, as she was very genuine with her answers and also because she changed majors as well and seems to be doing what she is passionate about. "
$ws.Range("G34").Value = "Tired with the overflow of info in the later part but the drawings for KEP are very interesting and well done!"

# --- Row 35 (15:00 - 16:00) ---
$ws.Rows.Item(35).RowHeight = 122.45
$ws.Range("A35").Value = 43883
$ws.Range("B35").Value = "15:00 - 16:00"
$ws.Range("C35").Value = "Team"
$ws.Range("D35").Value = "Discuss architecture, run over the assignment 2 resubmission"
$ws.Range("E35").Value = "We managed to look at the folder structures and settle on a MVC like pattern for as-implemented. As-intended architecture, we looked at closed pull requests/issues and were able to see the core developers take on it. Finished rewriting the Assignment 2 and submitted"
$ws.Range("F35").Value = "We noticed a lot of interdependencies among different components of the core structures in the software. This violates a true MVC design. We were happy to roughly come up with the architectural design after the fruitful discussion"
$ws.Range("G35").Value = "Challenging discussion"

# --- Row 36 (21:00 - 23:00) ---
$ws.Rows.Item(36).RowHeight = 188.55
$ws.Range("A36").Value = 43886
$ws.Range("B36").Value = "21:00 - 23:00"
$ws.Range("C36").Value = "Team"
$ws.Range("D36").Value = "Finish deciding on interesting pull requests and issues"
$ws.Range("E36").Value = "We managed to settle on five interesting pull requests and issues"
$ws.Range("F36").Value = "While reading through the various pull requests and issues, it was funny how the conversations are among the developers and the submitters. We also understand how not to write commit messages or how bland certain statements are. Also the feedback from the developers helps us in understanding the thought process they go through. We noticed a checklist for the submission which was interesting and the first time I had seen such a thing on GitHub.  "
$ws.Range("G36").Value = "Interesting!"

# --- Row 37 (21:00 - 22:00) ---
$ws.Rows.Item(37).RowHeight = 116.7
$ws.Range("A37").Value = 43887
$ws.Range("B37").Value = "21:00 - 22:00"
$ws.Range("C37").Value = "Team"
$ws.Range("D37").Value = "Finish write up on Social Context"
$ws.Range("E37").Value = "We managed to find resources that can aid us with developing the social context of the system"
$ws.Range("F37").Value = "We looked up the pulse in GitHub page and we were also surprised at how active the developers are, like merging within two days! Glad that we chose an active project and hopefully we can contribute effectively"
$ws.Range("G37").Value = "Satisfied!"

# --- Row 38 (8:00 - 12:00) ---
$ws.Rows.Item(38).RowHeight = 182.7
$ws.Range("A38").Value = 43888
$ws.Range("B38").Value = "8:00 - 12:00"
$ws.Range("C38").Value = "None"
$ws.Range("D38").Value = "Finish the entire write up, combine the individual parts we organized"
$ws.Range("E38").Value = "We managed to finish the write up, didn’t really think it would take this long as we had already discussed the key essence and architectural patterns"
$ws.Range("F38").Value = "There were a lot of confusion with regards to adhering to the architectural because we had interdependencies in a pure MVC. Hope we addressed these carefully in the report for Kajo to understand.  Social context section was a lot of work and effort. Nevertheless, we feel more confident about making a pull request and know where to go find the details. Feels like we know the project very well conceptually and somewhat working wise."
$ws.Range("G38").Value = "Praise the documentation and blog of JabRef."

Write-Output "Diary entries for Lecture 7 added (rows 34-38)."
